# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 3
$ws.Range("F11").Value = -2
$ws.Range("F17").Value = -3
$ws.Range("F25").Value = -4
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = 8
$ws.Range("F31").Value = -4
